# Daily attendance processing - 2025-11-16 07:21:28
# Rotate the "Recorded By" (column G) list of names/emails for each row:
# move the first comma-separated entry to the end of the list, for any
# cell that contains more than one comma-separated entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
